$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''37.470.77'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '''2.069.50'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''232.10'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").Value = '''0.628'
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''57.30'
$ws.Range("E8").Value = '  -1.86%  '
$ws.Range("D9").Value = '''0.389'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("D10").Value = '''0.0778'
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("E11").Value = '  +1.89%  '
$ws.Range("D12").Value = '''14.84'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '''2.372.13'
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").Value = '''20.87'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '''0.765'
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("D16").Value = '''5.32'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = '''2.067.82'
$ws.Range("E17").Value = '  -1.13%  '
$ws.Range("D18").Value = '''37.344.97'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").Value = '''70.41'
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").Value = '''5.97'
$ws.Range("E20").Value = '  -3.50%  '
$ws.Range("D21").Value = '''0.0₃0829'
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").Value = '''227.96'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '''2.35'
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("D25").Value = '''2.37'
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("D26").Value = '''9.62'
$ws.Range("E26").Value = '  +6.82%  '
$ws.Range("D27").Value = '''169.97'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").Value = '''0.133'
$ws.Range("E28").Value = '  -3.83%  '
$ws.Range("D29").Value = '''19.48'
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("D31").Value = '''0.123'
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").Value = '''4.61'
$ws.Range("E32").Value = '  -1.36%  '
$ws.Range("D33").Value = '''0.0632'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = '''4.63'
$ws.Range("E34").Value = '  -0.36%  '
$ws.Range("D35").Value = '''2.47'
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").Value = '''3.30'
$ws.Range("E37").Value = '  -2.69%  '
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = '''5.27'
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("D40").Value = '''0.0230'
$ws.Range("E40").Value = '  +7.01%  '
$ws.Range("D41").Value = '''99.61'
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").Value = '''0.0951'
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("E44").Value = '  +3.50%  '
$ws.Range("D45").Value = '''1.463.77'
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("D46").Value = '''16.70'
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("E47").Value = '  -1.47%  '
$ws.Range("D48").Value = '''3.95'
$ws.Range("E48").Value = '  -5.85%  '
$ws.Range("D49").Value = '''7.21'
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("D50").Value = '''2.94'
$ws.Range("E50").Value = '  -2.28%  '
$ws.Range("D51").Value = '''2.256.64'
$ws.Range("E51").Value = '  -0.45%  '
